# This script reproduces a weekly data refresh of the "Acelga" sheet:
# a new price-report row for the week is inserted at row 96 (pushing the
# existing rows 96-209 down to 97-210), and the new row is populated with
# the latest reported values while the rest of the historical rows keep
# their original data (now shifted down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 96; this shifts rows 96:209 down to 97:210 and
# expands the sheet dimension from A1:R209 to A1:R210.
$ws.Rows("96:96").Insert()

# Populate the newly inserted row 96 with the new weekly record.
$ws.Cells.Item(96, 1).Value  = 3
$ws.Cells.Item(96, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(96, 3).Value  = "Coquimbo"
$ws.Cells.Item(96, 4).Value  = 44483
$ws.Cells.Item(96, 5).Value  = 5
$ws.Cells.Item(96, 6).Value  = 100112009
$ws.Cells.Item(96, 7).Value  = "Acelga"
$ws.Cells.Item(96, 8).Value  = "Sin especificar"
$ws.Cells.Item(96, 9).Value  = "Primera"
$ws.Cells.Item(96, 10).Value = 250
$ws.Cells.Item(96, 11).Value = 2000
$ws.Cells.Item(96, 12).Value = 2200
$ws.Cells.Item(96, 13).Value = 2096
$ws.Cells.Item(96, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(96, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(96, 16).Value = 349
$ws.Cells.Item(96, 17).Value = 6
$ws.Cells.Item(96, 18).Value = "Hortaliza"
